# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting a refreshed data export (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet - column F updates by row
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2 = 1056
    3 = 660
    4 = 1464
    6 = 3212
    8 = 606
    9 = 2161
    10 = 462
    11 = 391
    12 = 230
    13 = 120
    14 = 277
    16 = 1050
    17 = 421
    20 = 4287
    21 = 1257
    22 = 3298
    24 = 135
    25 = 3157
    26 = 4785
    29 = 529
    30 = 3112
    31 = 320
    35 = 570
    36 = 1128
    37 = 0
    39 = 1280
    40 = 818
    42 = 763
    45 = 264
    47 = 115
    48 = 360
    49 = 3694
}
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# 本地生活 (Local life) sheet - column F updates by row
$ws3 = $wb.Worksheets.Item("本地生活")
$updates3 = @{
    2 = 1987
}
foreach ($row in $updates3.Keys) {
    $ws3.Cells.Item($row, 6).Value = $updates3[$row]
}

# 全部类型 (All types) sheet - column F updates by row
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2 = 1987
    3 = 660
    4 = 1464
    5 = 3212
    7 = 606
    9 = 2161
    10 = 462
    11 = 391
    12 = 230
    14 = 120
    15 = 277
    16 = 1050
    17 = 421
    19 = 4287
    21 = 1257
    23 = 3298
    24 = 3157
    25 = 4785
    28 = 3112
    29 = 320
    33 = 572
    34 = 1128
    37 = 1280
    39 = 818
    44 = 264
    47 = 115
    48 = 360
    49 = 3694
}
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
